# Add two new "relationship" sheets (UserAppliances, UserFuels) to the
# combined-import-template workbook, matching the Appliances/Fuels/Users
# sheets already present.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# UserAppliances
# ---------------------------------------------------------------------
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$ws1 = $wb.Worksheets.Add($null, $lastSheet)
$ws1.Name = "UserAppliances"

$ws1.Columns.Item(1).ColumnWidth = 15.83203125
$ws1.Columns.Item(2).ColumnWidth = 15.83203125
$ws1.Columns.Item(3).ColumnWidth = 15.83203125
$ws1.Columns.Item(4).ColumnWidth = 12.83203125
$ws1.Columns.Item(5).ColumnWidth = 40.83203125

$ws1.Range("A1").Value = "userId"
$ws1.Range("B1").Value = "applianceId"
$ws1.Range("C1").Value = "assignedDate"
$ws1.Range("D1").Value = "status"
$ws1.Range("E1").Value = "notes"

$ws1.Range("A2").Value = "USER001"
$ws1.Range("B2").Value = "APP001"
$ws1.Range("C2").Value = "15/01/2025"
$ws1.Range("D2").Value = "active"
$ws1.Range("E2").Value = "Primary heating appliance"

# ---------------------------------------------------------------------
# UserFuels
# ---------------------------------------------------------------------
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "UserFuels"

$ws2.Columns.Item(1).ColumnWidth = 15.83203125
$ws2.Columns.Item(2).ColumnWidth = 15.83203125
$ws2.Columns.Item(3).ColumnWidth = 15.83203125
$ws2.Columns.Item(4).ColumnWidth = 12.83203125
$ws2.Columns.Item(5).ColumnWidth = 40.83203125

$ws2.Range("A1").Value = "userId"
$ws2.Range("B1").Value = "fuelId"
$ws2.Range("C1").Value = "assignedDate"
$ws2.Range("D1").Value = "status"
$ws2.Range("E1").Value = "notes"

$ws2.Range("A2").Value = "USER001"
$ws2.Range("B2").Value = "FUEL001"
$ws2.Range("C2").Value = "15/01/2025"
$ws2.Range("D2").Value = "active"
$ws2.Range("E2").Value = "Preferred fuel type"
